# MonitoringShift.xlsx - reassign the Mar2017 on-call roster (columns B & C)
# to the updated names per the commit. Dates in column A, and the Feb2017
# sheet, are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mar2017")

$bNames = @("Arnel", "Junsat", "Carlo", "Leo", "Meryll", "Cath", "Brain", "Jun", "TinB", "Rodney", "Oscar", "Ivy", "Earl", "John", "Kennex", "Biboy", "Kate", "Marj", "Zhey", "Prado", "Mart", "Sky", "Kevin", "Roy", "Nathan", "Anj", "Reyn", "Jec", "Morgan", "Zhey", "Prado", "Junsat", "Ivy", "John", "Kate", "Kevin", "Sky", "Arnel", "Meryll", "Kennex", "TinB", "Cath", "Roy", "Morgan", "Mart", "Earl", "Carlo", "Marj", "Nathan", "Brain", "Anj", "Biboy", "Leo", "Carlo", "Prado", "Junsat", "Leo", "Arnel", "Meryll", "Brain", "Kevin", "Carlo")

$cNames = @("Cath", "Edch", "Meryll", "Kennex", "Lem", "John", "Roy", "Amy", "Morgan", "Leo", "Anj", "Ardeth", "Zhey", "Momay", "TinB", "Sky", "Pau", "Nathan", "Mart", "Kate", "Daisy", "Anne", "Ivy", "Harry", "TinC", "Claud", "Earl", "Biboy", "Anne", "Pati", "Mikee", "Eunice", "Marj", "Daisy", "Mikee", "Claud", "Jec", "Pati", "Rodney", "Ardeth", "Eunice", "Reyn", "Harry", "Pau", "TinC", "Pati", "Edch", "Harry", "Eunice", "Claud", "Momay", "Edch", "Jun", "Lem", "Amy", "Pau", "Mikee", "Oscar", "Anne", "Lem", "Momay", "Anne")

for ($i = 0; $i -lt $bNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bNames[$i]
    $ws.Cells.Item($row, 3).Value = $cNames[$i]
}
